{"js": "// Update the title date and the 25 answer cells (5 rows x 5 cols) of the\n// single answers table, per the commit's regenerated output.\n//\n// Row 0 of the table keeps its 5-cell shape in the final OOXML (tblGrid is\n// untouched), so every change is expressed as an in-place cell value write\n// addressed by (row, col) - never a global text search/replace - because\n// several old values coincide with other cells' new values (e.g. the old\n// row0/col1 text \"35\u00f79=3, 8\" becomes the new row0/col0 text, and the old\n// row16/col1 text \"76\u00f79=8, 4\" becomes the new row0/col0 text), which would\n// corrupt a sequential find/replace pass.\n\nconst title = context.document.body.paragraphs.getFirst();\ntitle.load('text');\n\nconst tables = context.document.body.tables;\ntables.load('items');\nawait context.sync();\n\n// --- Title date line ---\ntitle.insertText('2024-07-17 Wednesday', Word.InsertLocation.replace);\n\n// --- Answers table ---\nconst table = tables.items[0];\n\n// rowIndex -> new texts for the 5 columns, in document order\nconst newRowValues = {\n  0: ['35\u00f79=3, 8', '68\u00f76=11, 2', '78\u00f77=11, 1', '67\u00f79=7, 4', '53\u00f75=10, 3'],\n  4: ['90\u00f72=45, 0', '92\u00f74=23, 0', '54\u00f75=10, 4', '18\u00f77=2, 4', '77\u00f78=9, 5'],\n  8: ['76\u00f79=8, 4', '90\u00f75=18, 0', '24\u00f74=6, 0', '50\u00f77=7, 1', '15\u00f78=1, 7'],\n  12: ['34\u00f76=5, 4', '52\u00f75=10, 2', '42\u00f77=6, 0', '12\u00f72=6, 0', '36\u00f74=9, 0'],\n  16: ['97\u00f73=32, 1', '57\u00f72=28, 1', '83\u00f74=20, 3', '64\u00f76=10, 4', '95\u00f76=15, 5'],\n};\n\nfor (const rowIndex of Object.keys(newRowValues)) {\n  const r = Number(rowIndex);\n  const values = newRowValues[rowIndex];\n  for (let c = 0; c < values.length; c++) {\n    table.getCell(r, c).value = values[c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the title date and the 25 answer cells (5 rows x 5 cols) of the\n# single answers table, per the commit's regenerated output.\n#\n# Row 1 of the table (1-based COM indexing) keeps its 5-cell shape in the\n# final OOXML (the tblGrid is untouched), so every change is expressed as an\n# in-place cell text write addressed by (row, col) - never a blanket\n# Find/Replace across the whole document - because several old values\n# coincide with other cells' new values (e.g. the old row1/col2 text\n# \"35\u00f79=3, 8\" becomes the new row1/col1 text, and the old row17/col2 text\n# \"76\u00f79=8, 4\" becomes the new row1/col1 text), which would corrupt a\n# sequential find/replace pass.\n\n$d = $word.ActiveDocument\n\n# --- Title date line ---\n$d.Paragraphs(1).Range.Text = \"2024-07-17 Wednesday\"\n\n# --- Answers table ---\n$t = $d.Tables(1)\n\n# COM is 1-based: table row r (1-based) holds data for rows 1, 5, 9, 13, 17.\n$newRowValues = @{\n    1  = @(\"35\u00f79=3, 8\", \"68\u00f76=11, 2\", \"78\u00f77=11, 1\", \"67\u00f79=7, 4\", \"53\u00f75=10, 3\")\n    5  = @(\"90\u00f72=45, 0\", \"92\u00f74=23, 0\", \"54\u00f75=10, 4\", \"18\u00f77=2, 4\", \"77\u00f78=9, 5\")\n    9  = @(\"76\u00f79=8, 4\", \"90\u00f75=18, 0\", \"24\u00f74=6, 0\", \"50\u00f77=7, 1\", \"15\u00f78=1, 7\")\n    13 = @(\"34\u00f76=5, 4\", \"52\u00f75=10, 2\", \"42\u00f77=6, 0\", \"12\u00f72=6, 0\", \"36\u00f74=9, 0\")\n    17 = @(\"97\u00f73=32, 1\", \"57\u00f72=28, 1\", \"83\u00f74=20, 3\", \"64\u00f76=10, 4\", \"95\u00f76=15, 5\")\n}\n\nforeach ($rowIndex in $newRowValues.Keys) {\n    $values = $newRowValues[$rowIndex]\n    for ($c = 0; $c -lt $values.Length; $c++) {\n        $t.Cell($rowIndex, $c + 1).Range.Text = $values[$c]\n    }\n}\n"}
